$wb = $excel.ActiveWorkbook

$terms = $wb.Worksheets.Item("terms")
$terms.Range("B701").Value = "({float}|unknown)"
$terms.Range("B702").Value = "({text}{1,3}/{text}{1,3}|unknown)"

# Remove the "treatment", "race", "fitzpatrick", "additional data processing
# protocol", "labnotes as pdf" term rows (and the blank separator row after
# them) from the terms sheet.
$terms.Range("A735:A740").EntireRow.Delete() | Out-Null

$observationUnit = $wb.Worksheets.Item("ObservationUnit")
# Remove the "treatment", "race", "fitzpatrick" rows from the patient package.
$observationUnit.Range("A45:A47").EntireRow.Delete() | Out-Null

$assay = $wb.Worksheets.Item("Assay")
# Remove "additional data processing protocol" and "labnotes as pdf" rows
# from the metabolomics_T assay package.
$assay.Range("A204:A205").EntireRow.Delete() | Out-Null

